# Sync automático del tracker (cada 3h)
#
# Appends the newest closed/opened match rows (59-63) to the bottom of
# Sheet1, mirroring how the existing rows are laid out:
#   A event_id | B fecha | C jugador_A | D jugador_B | E pronostico | F cuota
# "resultado" (G) / "profit" (H) are intentionally left blank — those two
# columns only get filled in once a match result is known, same as the most
# recent rows already sitting at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# event_id/fecha look numeric/date-like to Excel's auto-detection, but the
# tracker stores them as plain text (matching every other row in the
# sheet), so force a text format right before writing them, then drop the
# format override again so the new rows don't end up with a stray style.
function Set-TextValue {
    param($Sheet, $Address, $Text)

    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

$newRows = @(
    @{ Row = 59; EventId = "14339497"; Fecha = "2025-08-05"; JugadorA = "Justin Engel";      JugadorB = "Alejandro Moro Canas"; Pronostico = "Gana Justin Engel";         Cuota = 3 },
    @{ Row = 60; EventId = "14349616"; Fecha = "2025-08-05"; JugadorA = "Abdullah Shelbayh"; JugadorB = "Ugo Blanchet";          Pronostico = "Gana Abdullah Shelbayh";   Cuota = 2.75 },
    @{ Row = 61; EventId = "14349615"; Fecha = "2025-08-05"; JugadorA = "Rafael Jodar";      JugadorB = "Marc-Andrea Huesler";   Pronostico = "Gana Marc-Andrea Huesler"; Cuota = 1.91 },
    @{ Row = 62; EventId = "14349721"; Fecha = "2025-08-05"; JugadorA = "Max Alcala Gurri";  JugadorB = "Lorenzo Giustino";      Pronostico = "Gana Lorenzo Giustino";    Cuota = 2.1 },
    @{ Row = 63; EventId = "14349731"; Fecha = "2025-08-05"; JugadorA = "Lukas Neumayer";    JugadorB = "Mariano Kestelboim";    Pronostico = "Gana Mariano Kestelboim";  Cuota = 3.75 }
)

foreach ($row in $newRows) {
    $r = $row.Row

    Set-TextValue $ws "A$r" $row.EventId
    Set-TextValue $ws "B$r" $row.Fecha
    $ws.Range("C$r").Value = $row.JugadorA
    $ws.Range("D$r").Value = $row.JugadorB
    $ws.Range("E$r").Value = $row.Pronostico
    $ws.Range("F$r").Value = $row.Cuota
    # G$r (resultado) / H$r (profit) stay blank until the match is settled.
}
